# Apply the Natmi L1cam-Alcam LR-pair table update ("Natmi following Dr Hou advice"):
# expand the Sending/Target cluster cross-join from 3 clusters to the full 4x4 grid
# (ECs, FAPs, M2, sCs) and refresh the recomputed expression/specificity metrics
# for the L1cam -> Alcam ligand-receptor pair (rows 2-17, columns A:T).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "L1cam", "Alcam", "ECs", 3, 1, 18.04537966666667, 54.13613900000001, 0.6797959733292525, 0.6797959733292525, 2, 0.6666666666666666, 29.77446566666667, 89.323397, 0.4936227575398824, 0.4936227575398824, 537.2915373271316, 4835.623835944184, 0.3355627629192939, 0.3355627629192939),
    @("ECs", "L1cam", "Alcam", "FAPs", 3, 1, 18.04537966666667, 54.13613900000001, 0.6797959733292525, 0.6797959733292525, 3, 1, 0.5058753333333333, 1.517626, 0.00838676938175808, 0.00838676938175808, 9.128712454001557, 82.15841208601402, 0.005701292054960208, 0.005701292054960208),
    @("ECs", "L1cam", "Alcam", "M2", 3, 1, 18.04537966666667, 54.13613900000001, 0.6797959733292525, 0.6797959733292525, 3, 1, 29.12203633333333, 87.366109, 0.482806309304492, 0.482806309304492, 525.5182023014612, 4729.663820713151, 0.3282097849631512, 0.3282097849631512),
    @("ECs", "L1cam", "Alcam", "sCs", 3, 1, 18.04537966666667, 54.13613900000001, 0.6797959733292525, 0.6797959733292525, 3, 1, 0.9158823333333334, 2.747647, 0.0151841637738675, 0.0151841637738675, 16.52744443499256, 148.746999914933, 0.01032213339184703, 0.01032213339184703),
    @("FAPs", "L1cam", "Alcam", "ECs", 3, 1, 0.6001993333333334, 1.800598, 0.02261039099934159, 0.02261039099934159, 2, 0.6666666666666666, 29.77446566666667, 89.323397, 0.4936227575398824, 0.4936227575398824, 17.87061444348956, 160.835529991406, 0.01116100355414993, 0.01116100355414993),
    @("FAPs", "L1cam", "Alcam", "FAPs", 3, 1, 0.6001993333333334, 1.800598, 0.02261039099934159, 0.02261039099934159, 3, 1, 0.5058753333333333, 1.517626, 0.00838676938175808, 0.00838676938175808, 0.3036260378164445, 2.732634340348, 0.0001896281349428566, 0.0001896281349428566),
    @("FAPs", "L1cam", "Alcam", "M2", 3, 1, 0.6001993333333334, 1.800598, 0.02261039099934159, 0.02261039099934159, 3, 1, 29.12203633333333, 87.366109, 0.482806309304492, 0.482806309304492, 17.47902679257578, 157.311241133182, 0.01091643943032362, 0.01091643943032362),
    @("FAPs", "L1cam", "Alcam", "sCs", 3, 1, 0.6001993333333334, 1.800598, 0.02261039099934159, 0.02261039099934159, 3, 1, 0.9158823333333334, 2.747647, 0.0151841637738675, 0.0151841637738675, 0.5497119658784445, 4.947407692906001, 0.0003433198799251824, 0.0003433198799251825),
    @("M2", "L1cam", "Alcam", "ECs", 3, 1, 4.690054333333333, 14.070163, 0.1766812397072912, 0.1766812397072912, 2, 0.6666666666666666, 29.77446566666667, 89.323397, 0.4936227575398824, 0.4936227575398824, 139.6438617226346, 1256.794755503711, 0.08721388074987803, 0.08721388074987804),
    @("M2", "L1cam", "Alcam", "FAPs", 3, 1, 4.690054333333333, 14.070163, 0.1766812397072912, 0.1766812397072912, 3, 1, 0.5058753333333333, 1.517626, 0.00838676938175808, 0.00838676938175808, 2.372582799226445, 21.353245193038, 0.00148178481150817, 0.00148178481150817),
    @("M2", "L1cam", "Alcam", "M2", 3, 1, 4.690054333333333, 14.070163, 0.1766812397072912, 0.1766812397072912, 3, 1, 29.12203633333333, 87.366109, 0.482806309304492, 0.482806309304492, 136.5839327006408, 1229.255394305767, 0.0853028172664195, 0.08530281726641951),
    @("M2", "L1cam", "Alcam", "sCs", 3, 1, 4.690054333333333, 14.070163, 0.1766812397072912, 0.1766812397072912, 3, 1, 0.9158823333333334, 2.747647, 0.0151841637738675, 0.0151841637738675, 4.295537906273444, 38.65984115646101, 0.002682756879485451, 0.002682756879485452),
    @("sCs", "L1cam", "Alcam", "ECs", 3, 1, 3.209654333333333, 9.628962999999999, 0.1209123959641148, 0.1209123959641148, 2, 0.6666666666666666, 29.77446566666667, 89.323397, 0.4936227575398824, 0.4936227575398824, 95.56574274970122, 860.0916847473109, 0.05968511031656049, 0.05968511031656049),
    @("sCs", "L1cam", "Alcam", "FAPs", 3, 1, 3.209654333333333, 9.628962999999999, 0.1209123959641148, 0.1209123959641148, 3, 1, 0.5058753333333333, 1.517626, 0.00838676938175808, 0.00838676938175808, 1.623684955759778, 14.613164601838, 0.001014064380346847, 0.001014064380346847),
    @("sCs", "L1cam", "Alcam", "M2", 3, 1, 3.209654333333333, 9.628962999999999, 0.1209123959641148, 0.1209123959641148, 3, 1, 29.12203633333333, 87.366109, 0.482806309304492, 0.482806309304492, 93.47167011277409, 841.2450310149668, 0.05837726764459761, 0.05837726764459761),
    @("sCs", "L1cam", "Alcam", "sCs", 3, 1, 3.209654333333333, 9.628962999999999, 0.1209123959641148, 0.1209123959641148, 3, 1, 0.9158823333333334, 2.747647, 0.0151841637738675, 0.0151841637738675, 2.939665700006778, 26.456991300061, 0.001835953622609835, 0.001835953622609835)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowVals = $rows[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value2 = $rowVals[$c]
    }
}

Write-Host "Updated rows 2-17 (A:T) on $($ws.Name)"
